$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 58338292
$ws.Range("I62").Value = 25005328
$ws.Range("J62").Value = 125004220
$ws.Range("K62").Value = 25005328
$ws.Range("L62").Value = 125004220
$ws.Range("M62").Value = -25004704
$ws.Range("N62").Value = -125005468

$ws.Range("H65").Value = 58338292
$ws.Range("I65").Value = 25005328
$ws.Range("J65").Value = 125004220
$ws.Range("K65").Value = 125026640
$ws.Range("L65").Value = 625021100
$ws.Range("M65").Value = -125023520
$ws.Range("N65").Value = -625027340

$ws.Range("H107").Value = 1505.9286
$ws.Range("I107").Value = 1908.3
$ws.Range("K107").Value = 1908.3
$ws.Range("M107").Value = 11.70000000000005

$ws.Range("H137").Value = 20284364
$ws.Range("I137").Value = 4033133.8
$ws.Range("J137").Value = 104249050
$ws.Range("K137").Value = 12099401.4
$ws.Range("L137").Value = 312747150
$ws.Range("M137").Value = -12096851.4
$ws.Range("N137").Value = -312752250


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 295118.78
$ws.Range("I45").Value = 556421.5600000001
$ws.Range("J45").Value = 1153.1875
$ws.Range("K45").Value = 556421.5600000001
$ws.Range("L45").Value = 1153.1875
$ws.Range("M45").Value = -556044.5600000001
$ws.Range("N45").Value = -1907.1875

$ws.Range("H61").Value = 4369967.5
$ws.Range("I61").Value = 1984987.1
$ws.Range("J61").Value = 29412264
$ws.Range("K61").Value = 1984987.1
$ws.Range("L61").Value = 29412264
$ws.Range("M61").Value = -1984775.1
$ws.Range("N61").Value = -29412688

$ws.Range("H74").Value = 57321044
$ws.Range("I74").Value = 48299980
$ws.Range("J74").Value = 88894776
$ws.Range("K74").Value = 48299980
$ws.Range("L74").Value = 88894776
$ws.Range("M74").Value = -48299106
$ws.Range("N74").Value = -88896524

$ws.Range("H77").Value = 57321044
$ws.Range("I77").Value = 48299980
$ws.Range("J77").Value = 88894776
$ws.Range("K77").Value = 241499900
$ws.Range("L77").Value = 444473880
$ws.Range("M77").Value = -241495532
$ws.Range("N77").Value = -444482616

$ws.Range("H88").Value = 5490
$ws.Range("I88").Value = 2350
$ws.Range("J88").Value = 7583.3335
$ws.Range("K88").Value = 2350
$ws.Range("L88").Value = 7583.3335
$ws.Range("M88").Value = -1944
$ws.Range("N88").Value = -8395.333500000001

$ws.Range("H91").Value = 5490
$ws.Range("I91").Value = 2350
$ws.Range("J91").Value = 7583.3335
$ws.Range("K91").Value = 2350
$ws.Range("L91").Value = 7583.3335
$ws.Range("M91").Value = -946
$ws.Range("N91").Value = -10391.3335

$ws.Range("H110").Value = 1113.8334
$ws.Range("I110").Value = 1170.75
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 1170.75
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 874.25
$ws.Range("N110").Value = -5090

$ws.Range("H132").Value = 11788780
$ws.Range("I132").Value = 13338155
$ws.Range("J132").Value = 6946982
$ws.Range("K132").Value = 40014465
$ws.Range("L132").Value = 20840946
$ws.Range("M132").Value = -40011935
$ws.Range("N132").Value = -20846006

$ws.Range("H136").Value = 4369967.5
$ws.Range("I136").Value = 1984987.1
$ws.Range("J136").Value = 29412264
$ws.Range("K136").Value = 5954961.300000001
$ws.Range("L136").Value = 88236792
$ws.Range("M136").Value = -5952411.300000001
$ws.Range("N136").Value = -88241892


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1926.36
$ws.Range("I86").Value = 1939.2766
$ws.Range("J86").Value = 1724
$ws.Range("K86").Value = 1939.2766
$ws.Range("L86").Value = 1724
$ws.Range("M86").Value = -816.2765999999999
$ws.Range("N86").Value = -3970

$ws.Range("H89").Value = 1926.36
$ws.Range("I89").Value = 1939.2766
$ws.Range("J89").Value = 1724
$ws.Range("K89").Value = 9696.383
$ws.Range("L89").Value = 8620
$ws.Range("M89").Value = -4080.383
$ws.Range("N89").Value = -19852

$ws.Range("H105").Value = 1790.2941
$ws.Range("I105").Value = 1762.3334
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1762.3334
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -15.33339999999998
$ws.Range("N105").Value = -5494

$ws.Range("H134").Value = 26787720
$ws.Range("I134").Value = 35715730
$ws.Range("J134").Value = 5955699.5
$ws.Range("K134").Value = 107147190
$ws.Range("L134").Value = 17867098.5
$ws.Range("M134").Value = -107144655
$ws.Range("N134").Value = -17872168.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2543894.5
$ws.Range("I31").Value = 1303589.5
$ws.Range("J31").Value = 6953868
$ws.Range("K31").Value = 1303589.5
$ws.Range("L31").Value = 6953868
$ws.Range("M31").Value = -1303294.5
$ws.Range("N31").Value = -6954458

$ws.Range("H34").Value = 2543894.5
$ws.Range("I34").Value = 1303589.5
$ws.Range("J34").Value = 6953868
$ws.Range("K34").Value = 1303589.5
$ws.Range("L34").Value = 6953868
$ws.Range("M34").Value = -1303387.5
$ws.Range("N34").Value = -6954272

$ws.Range("H58").Value = 4778956
$ws.Range("I58").Value = 2752476.5
$ws.Range("J58").Value = 11365014
$ws.Range("K58").Value = 2752476.5
$ws.Range("L58").Value = 11365014
$ws.Range("M58").Value = -2752273.5
$ws.Range("N58").Value = -11365420

$ws.Range("H62").Value = 3099.1177
$ws.Range("I62").Value = 2640
$ws.Range("J62").Value = 3755
$ws.Range("K62").Value = 2640
$ws.Range("L62").Value = 3755
$ws.Range("M62").Value = -2016
$ws.Range("N62").Value = -5003

$ws.Range("H65").Value = 3099.1177
$ws.Range("I65").Value = 2640
$ws.Range("J65").Value = 3755
$ws.Range("K65").Value = 13200
$ws.Range("L65").Value = 18775
$ws.Range("M65").Value = -10080
$ws.Range("N65").Value = -25015

$ws.Range("H132").Value = 1472905
$ws.Range("I132").Value = 2084990
$ws.Range("J132").Value = 3901.1
$ws.Range("K132").Value = 6254970
$ws.Range("L132").Value = 11703.3
$ws.Range("M132").Value = -6252440
$ws.Range("N132").Value = -16763.3

$ws.Range("H134").Value = 1217315.4
$ws.Range("I134").Value = 5606.2173
$ws.Range("J134").Value = 4004246.5
$ws.Range("K134").Value = 16818.6519
$ws.Range("L134").Value = 12012739.5
$ws.Range("M134").Value = -14283.6519
$ws.Range("N134").Value = -12017809.5

$ws.Range("H136").Value = 4778956
$ws.Range("I136").Value = 2752476.5
$ws.Range("J136").Value = 11365014
$ws.Range("K136").Value = 8257429.5
$ws.Range("L136").Value = 34095042
$ws.Range("M136").Value = -8254879.5
$ws.Range("N136").Value = -34100142


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1660.5128
$ws.Range("I113").Value = 1124.1875
$ws.Range("J113").Value = 2033.6086
$ws.Range("K113").Value = 3372.5625
$ws.Range("L113").Value = 6100.825800000001
$ws.Range("M113").Value = -1202.5625
$ws.Range("N113").Value = -10440.8258


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 23813010
$ws.Range("J122").Value = 55557724
$ws.Range("L122").Value = 166673172
$ws.Range("N122").Value = -166678072

$ws.Range("H132").Value = 26042650
$ws.Range("I132").Value = 35375036
$ws.Range("J132").Value = 15154864
$ws.Range("K132").Value = 106125108
$ws.Range("L132").Value = 45464592
$ws.Range("M132").Value = -106122578
$ws.Range("N132").Value = -45469652


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1599.5714
$ws.Range("I7").Value = 1439.6
$ws.Range("K7").Value = 1439.6
$ws.Range("M7").Value = -1327.6

$ws.Range("H122").Value = 10539661
$ws.Range("I122").Value = 1522492.2
$ws.Range("J122").Value = 28574000
$ws.Range("K122").Value = 4567476.6
$ws.Range("L122").Value = 85722000
$ws.Range("M122").Value = -4565026.6
$ws.Range("N122").Value = -85726900

$ws.Range("H126").Value = 1599.5714
$ws.Range("I126").Value = 1439.6
$ws.Range("K126").Value = 4318.799999999999
$ws.Range("M126").Value = -1848.799999999999

$ws.Range("H132").Value = 3179102.5
$ws.Range("I132").Value = 5129927
$ws.Range("J132").Value = 9012.25
$ws.Range("K132").Value = 15389781
$ws.Range("L132").Value = 27036.75
$ws.Range("M132").Value = -15387251
$ws.Range("N132").Value = -32096.75

$ws.Range("H136").Value = 13624540
$ws.Range("I136").Value = 25884312
$ws.Range("J136").Value = 2571.6667
$ws.Range("K136").Value = 77652936
$ws.Range("L136").Value = 7715.000100000001
$ws.Range("M136").Value = -77650386
$ws.Range("N136").Value = -12815.0001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 585141.6
$ws.Range("I132").Value = 2218.65
$ws.Range("K132").Value = 6655.950000000001
$ws.Range("M132").Value = -4125.950000000001

$ws.Range("H136").Value = 9044.071
$ws.Range("I136").Value = 6491.95
$ws.Range("J136").Value = 15424.375
$ws.Range("K136").Value = 19475.85
$ws.Range("L136").Value = 46273.125
$ws.Range("M136").Value = -16925.85
$ws.Range("N136").Value = -51373.125

